$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- At-bat #1 (rows 10-17) ---
$ws.Range("F10").Value = "FB"
$ws.Range("G10").Value = "Swing"
$ws.Range("H10").Value = "In Play"
$ws.Range("M10").Value = "76.77 MPH"

$ws.Range("M12").Value = "68.99°"

$ws.Range("J17").Value = "CH,CB,FB"

# --- At-bat #2 (rows 19-26) ---
$ws.Range("F19").Value = "FB"
$ws.Range("G19").Value = "Take"
$ws.Range("H19").Value = "Strike"

$ws.Range("F20").Value = "CB"
$ws.Range("G20").Value = "Take"
$ws.Range("H20").Value = "Ball"

$ws.Range("F21").Value = "CB"
$ws.Range("G21").Value = "Take"
$ws.Range("H21").Value = "Strike"
$ws.Range("M21").Value = $null

$ws.Range("F22").Value = "CH"
$ws.Range("G22").Value = "Take"
$ws.Range("H22").Value = "Ball"

$ws.Range("F23").Value = "FB"
$ws.Range("G23").Value = "Swing"
$ws.Range("H23").Value = "Foul"

$ws.Range("F24").Value = "CB"
$ws.Range("G24").Value = "Take"
$ws.Range("H24").Value = "Ball"
$ws.Range("M24").Value = "Walk"

$ws.Range("F25").Value = "FB"
$ws.Range("G25").Value = "Take"
$ws.Range("H25").Value = "Ball"

$ws.Range("J26").Value = "CH,CB,FB"

# --- At-bat #3 (rows 28-35) ---
$ws.Range("F28").Value = "CH"
$ws.Range("G28").Value = "Swing"
$ws.Range("H28").Value = "In Play"
$ws.Range("M28").Value = "80.67 MPH"

$ws.Range("M30").Value = "58.13°"

$ws.Range("J35").Value = "CH,CB,FB,SL"

# --- At-bat #4 (rows 37-44) ---
$ws.Range("F37").Value = "CH"
$ws.Range("G37").Value = "Take"
$ws.Range("H37").Value = "Strike"
$ws.Range("M37").Value = "76.77 MPH"

$ws.Range("F38").Value = "CB"
$ws.Range("G38").Value = "Take"
$ws.Range("H38").Value = "Strike"

$ws.Range("F39").Value = "CB"
$ws.Range("G39").Value = "Swing"
$ws.Range("H39").Value = "Foul"
$ws.Range("M39").Value = "21.94°"

$ws.Range("F40").Value = "CH"
$ws.Range("G40").Value = "Swing"
$ws.Range("H40").Value = "In Play"

$ws.Range("J44").Value = "CH,CB,FB,SL"
